$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text semantics on the Price/Volume columns so numeric-looking strings
# (e.g. "1.00", "311.92") are preserved as text instead of being parsed as numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '46.083.38'
$ws.Range('E2').Value = '  -0.48%  '
$ws.Range('D3').Value = '2.596.78'
$ws.Range('E3').Value = '  -0.31%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '311.92'
$ws.Range('E5').Value = '  +1.43%  '
$ws.Range('D6').Value = '99.00'
$ws.Range('E6').Value = '  -1.38%  '
$ws.Range('E7').Value = '  -0.19%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').Value = '0.584'
$ws.Range('E9').Value = '  +1.87%  '
$ws.Range('E10').Value = '  +0.03%  '
$ws.Range('D11').Value = '54.19'
$ws.Range('E11').Value = '  -1.87%  '
$ws.Range('E12').Value = '  +0.58%  '
$ws.Range('D13').Value = '8.14'
$ws.Range('E13').Value = '  -0.12%  '
$ws.Range('D14').Value = '2.995.10'
$ws.Range('E14').Value = '  -0.36%  '
$ws.Range('E15').Value = '  +1.45%  '
$ws.Range('D16').Value = '2.591.46'
$ws.Range('E16').Value = '  -1.30%  '
$ws.Range('E17').Value = '  +2.06%  '
$ws.Range('D18').Value = '14.87'
$ws.Range('E18').Value = '  -0.18%  '
$ws.Range('D19').Value = '46.209.60'
$ws.Range('E19').Value = '  -0.54%  '
$ws.Range('E20').Value = '  +1.29%  '
$ws.Range('D21').Value = '6.74'
$ws.Range('E21').Value = '  +0.92%  '
$ws.Range('D22').Value = '12.82'
$ws.Range('E22').Value = '  -3.24%  '
$ws.Range('D23').Value = '296.53'
$ws.Range('E23').Value = '  +16.06%  '
$ws.Range('D24').Value = '73.07'
$ws.Range('E24').Value = '  +2.76%  '
$ws.Range('D25').Value = '3.06'
$ws.Range('E25').Value = '  +1.52%  '
$ws.Range('E26').Value = '  -0.28%  '
$ws.Range('D27').Value = '29.55'
$ws.Range('E27').Value = '  +4.68%  '
$ws.Range('E28').Value = '  -0.02%  '
$ws.Range('E29').Value = '  +1.03%  '
$ws.Range('D30').Value = '10.83'
$ws.Range('E30').Value = '  +3.64%  '
$ws.Range('D31').Value = '38.81'
$ws.Range('E31').Value = '  -2.99%  '
$ws.Range('E32').Value = '  -2.43%  '
$ws.Range('D33').Value = '6.21'
$ws.Range('E33').Value = '  +0.87%  '
$ws.Range('D34').Value = '3.60'
$ws.Range('E34').Value = '  -3.19%  '
$ws.Range('D35').Value = '155.93'
$ws.Range('E35').Value = '  +2.89%  '
$ws.Range('D36').Value = '0.0837'
$ws.Range('E36').Value = '  +0.31%  '
$ws.Range('E37').Value = '  -5.95%  '
$ws.Range('E38').Value = '  -5.24%  '
$ws.Range('E39').Value = '  +3.97%  '
$ws.Range('D40').Value = '0.124'
$ws.Range('E40').Value = '  +1.33%  '
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').Value = '0.0332'
$ws.Range('E41').Value = '  +2.95%  '
$ws.Range('B42').Value = 'Celestia'
$ws.Range('C42').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D42').Value = '15.71'
$ws.Range('E42').Value = '  +0.51%  '
$ws.Range('D43').Value = '3.59'
$ws.Range('E43').Value = '  -0.28%  '
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').Value = '3.98'
$ws.Range('E44').Value = '  -4.85%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = '21.07'
$ws.Range('E45').Value = '  +8.88%  '
$ws.Range('D46').Value = '2.109.19'
$ws.Range('E46').Value = '  +2.46%  '
$ws.Range('D47').Value = '98.56'
$ws.Range('E47').Value = '  +8.32%  '
$ws.Range('D48').Value = '0.998'
$ws.Range('E48').Value = '  -0.11%  '
$ws.Range('D49').Value = '9.53'
$ws.Range('E49').Value = '  +4.25%  '
$ws.Range('B50').Value = 'Aave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D50').Value = '108.71'
$ws.Range('E50').Value = '  -0.80%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').Value = '0.202'
$ws.Range('E51').Value = '  +0.99%  '

# Restore the default (unstyled) cell style now that the values are committed as text,
# matching the original workbook which has no explicit style on these cells.
$ws.Range("D2:E51").Style = "Normal"
